# verigen 에 vfunction 추가
# Adds a new "vfunction" worksheet documenting the `vfunction(name, func)` verilog macro,
# references it from the "summary" sheet, and leaves "vfunction" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) summary sheet: append a new row describing "vfunction"
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("summary")

$summary.Range("A3").Value = "vfunction"
$summary.Range("B3").Value = "function"
$summary.Range("B3").HorizontalAlignment = -4108   # xlCenter, matches B1/B2 style
$summary.Range("C3").Value = "verilog 내 사용 함수 선언"

# ------------------------------------------------------------------
# 2) add the new "vfunction" worksheet, placed after "_V" (last tab)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$vfunction = $wb.Worksheets.Add($null, $lastSheet)
$vfunction.Name = "vfunction"

$vfunction.Range("A1").Value = "구분"
$vfunction.Range("B1").Value = "표현식"

$vfunction.Range("A2").Value = "함수 원형"
$vfunction.Range("A3").Value = "반환값"
$vfunction.Range("A4").Value = "설명"

$vfunction.Range("A5").Value = "name"
$vfunction.Range("B2").Value = "function vfunction(name, func)"

$vfunction.Range("A6").Value = "func"
$vfunction.Range("B6").Value = "verilog에서 사용할 lua function"

$vfunction.Range("B5").Value = "verilog 내에서 사용할 함수 이름."
$vfunction.Range("B3").Value = "-"
$vfunction.Range("B4").Value = 'verilog 내에서 "$함수(...)" 로 lua 함수를 호출할 수 있습니다.'

# Column A holds the same "header" style as the other sheets (centered),
# and the header cell B1 matches it too.
$vfunction.Range("A1:A6").HorizontalAlignment = -4108   # xlCenter
$vfunction.Range("B1").HorizontalAlignment = -4108      # xlCenter (header cell)

# Description cell wraps, like the analogous cell on the "_V" sheet.
$vfunction.Range("B4").WrapText = $true

# Final selections, matching the end-of-edit-session cursor state.
$summary.Range("B3").Select()
$vfunction.Range("B5").Select()
$vfunction.Activate()
